$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "39.952.25"
Set-TextValue $ws.Range("E2") "  +0.50%  "

Set-TextValue $ws.Range("D3") "2.211.58"
Set-TextValue $ws.Range("E3") "  -0.37%  "

Set-TextValue $ws.Range("E4") "  +0.08%  "

Set-TextValue $ws.Range("D5") "288.65"
Set-TextValue $ws.Range("E5") "  -1.91%  "

Set-TextValue $ws.Range("D6") "87.50"
Set-TextValue $ws.Range("E6") "  +4.18%  "

Set-TextValue $ws.Range("D7") "0.515"
Set-TextValue $ws.Range("E7") "  -0.17%  "

Set-TextValue $ws.Range("E8") "  -0.06%  "

Set-TextValue $ws.Range("D9") "0.469"
Set-TextValue $ws.Range("E9") "  +0.64%  "

Set-TextValue $ws.Range("D10") "30.53"
Set-TextValue $ws.Range("E10") "  +2.83%  "

Set-TextValue $ws.Range("D11") "0.0777"
Set-TextValue $ws.Range("E11") "  -0.75%  "

Set-TextValue $ws.Range("E12") "  +2.53%  "

Set-TextValue $ws.Range("D13") "6.44"
Set-TextValue $ws.Range("E13") "  +2.40%  "

Set-TextValue $ws.Range("D14") "2.554.44"
Set-TextValue $ws.Range("E14") "  -0.33%  "

Set-TextValue $ws.Range("D15") "13.94"
Set-TextValue $ws.Range("E15") "  -1.24%  "

Set-TextValue $ws.Range("D16") "2.213.54"
Set-TextValue $ws.Range("E16") "  -0.33%  "

Set-TextValue $ws.Range("D17") "0.727"
Set-TextValue $ws.Range("E17") "  +1.06%  "

Set-TextValue $ws.Range("D18") "39.895.55"
Set-TextValue $ws.Range("E18") "  +0.63%  "

Set-TextValue $ws.Range("D19") "11.66"
Set-TextValue $ws.Range("E19") "  +11.79%  "

Set-TextValue $ws.Range("D20") "0.0₃0883"
Set-TextValue $ws.Range("E20") "  +0.19%  "

Set-TextValue $ws.Range("E21") "  +0.93%  "

Set-TextValue $ws.Range("D22") "65.42"
Set-TextValue $ws.Range("E22") "  +0.65%  "

Set-TextValue $ws.Range("D23") "235.02"
Set-TextValue $ws.Range("E23") "  +1.18%  "

Set-TextValue $ws.Range("D24") "0.999"
Set-TextValue $ws.Range("E24") "  -0.08%  "

Set-TextValue $ws.Range("D25") "2.44"
Set-TextValue $ws.Range("E25") "  +1.60%  "

Set-TextValue $ws.Range("E26") "  +0.79%  "

Set-TextValue $ws.Range("B27") "EthereumClassic"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D27") "22.55"
Set-TextValue $ws.Range("E27") "  -1.18%  "

Set-TextValue $ws.Range("B28") "Toncoin"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D28") "2.19"
Set-TextValue $ws.Range("E28") "  +0.26%  "

Set-TextValue $ws.Range("D29") "9.21"
Set-TextValue $ws.Range("E29") "  +0.36%  "

Set-TextValue $ws.Range("D30") "152.86"
Set-TextValue $ws.Range("E30") "  +2.06%  "

Set-TextValue $ws.Range("D31") "31.88"
Set-TextValue $ws.Range("E31") "  -0.94%  "

Set-TextValue $ws.Range("E32") "  -0.02%  "

Set-TextValue $ws.Range("D33") "4.95"
Set-TextValue $ws.Range("E33") "  +2.76%  "

Set-TextValue $ws.Range("D34") "0.0717"
Set-TextValue $ws.Range("E34") "  +2.00%  "

Set-TextValue $ws.Range("E35") "  +0.67%  "

Set-TextValue $ws.Range("E36") "  +5.80%  "

Set-TextValue $ws.Range("E37") "  +0.15%  "

Set-TextValue $ws.Range("B38") "Celestia"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D38") "15.77"
Set-TextValue $ws.Range("E38") "  -1.29%  "

Set-TextValue $ws.Range("B39") "Kaspa"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D39") "0.0991"
Set-TextValue $ws.Range("E39") "  +2.30%  "

Set-TextValue $ws.Range("D40") "1.70"
Set-TextValue $ws.Range("E40") "  +3.05%  "

Set-TextValue $ws.Range("D41") "3.86"
Set-TextValue $ws.Range("E41") "  +5.06%  "

Set-TextValue $ws.Range("D42") "2.102.19"
Set-TextValue $ws.Range("E42") "  +8.39%  "

Set-TextValue $ws.Range("E43") "  +2.10%  "

Set-TextValue $ws.Range("D44") "0.0267"
Set-TextValue $ws.Range("E44") "  +0.57%  "

Set-TextValue $ws.Range("D45") "9.96"
Set-TextValue $ws.Range("E45") "  +5.80%  "

Set-TextValue $ws.Range("D46") "17.54"
Set-TextValue $ws.Range("E46") "  +8.39%  "

Set-TextValue $ws.Range("D47") "2.66"
Set-TextValue $ws.Range("E47") "  +2.59%  "

Set-TextValue $ws.Range("D48") "2.429.81"
Set-TextValue $ws.Range("E48") "  -0.09%  "

Set-TextValue $ws.Range("D49") "1.45"
Set-TextValue $ws.Range("E49") "  +2.13%  "

Set-TextValue $ws.Range("D50") "69.30"
Set-TextValue $ws.Range("E50") "  -1.93%  "

Set-TextValue $ws.Range("D51") "88.46"
Set-TextValue $ws.Range("E51") "  -0.47%  "
